$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the top of this data block
# (row 679). Insert a new row there, which pushes the existing rows
# 679:784 down to 680:785 (dimension grows from A1:R784 to A1:R785),
# then populate the new row with the latest reading.
$ws.Rows.Item(679).Insert()

$ws.Cells.Item(679, 1).Value = 6
$ws.Cells.Item(679, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(679, 3).Value = "Metropolitana"
$ws.Cells.Item(679, 4).Value = 45218
$ws.Cells.Item(679, 5).Value = 13
$ws.Cells.Item(679, 6).Value = 100112039
$ws.Cells.Item(679, 7).Value = "Ciboulette"
$ws.Cells.Item(679, 8).Value = "Sin especificar"
$ws.Cells.Item(679, 9).Value = "Primera"
$ws.Cells.Item(679, 10).Value = 580
$ws.Cells.Item(679, 11).Value = 1200
$ws.Cells.Item(679, 12).Value = 1300
$ws.Cells.Item(679, 13).Value = 1257
$ws.Cells.Item(679, 14).Value = "`$/docena de atados"
$ws.Cells.Item(679, 15).Value = "Región Metropolitana"
$ws.Cells.Item(679, 16).Value = 419
$ws.Cells.Item(679, 17).Value = 3
$ws.Cells.Item(679, 18).Value = "Hortaliza"
